$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-37 of the FuelGroup/Year demand table: the fuel-group taxonomy was
# corrected (Fossil Gases + Fossil Liquids split out of the totals) and the
# Pass Aviation (column H) figures were filled in for every fuel group, for
# each of the three model years (2030/2040/2050). Rewriting the full A:K block
# per row (instead of only touching previously-populated cells) keeps the
# row order, the row-to-year mapping and the blank/non-blank pattern correct.
$rows = @(
    @("Fossil Gases", 2030, $null, $null, $null, 0.0003497013861488061, $null, $null, [double]"3.436751862763012e-05", $null, $null),
    @("Synthetic Liquids", 2030, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Biogenic Liquids", 2030, $null, $null, $null, 0.001636165311418535, [double]"5.882606390480593e-06", 0.0003132736250606, 0.0016392156151277, [double]"3.556733278696332e-05", 0.002177464365623002),
    @("Fossil Liquids", 2030, $null, $null, $null, 0.01617614087513679, [double]"4.142947176255379e-05", 0.002863152017891, 0.0102584426004798, 0.0002152891401867, 0.02132334971103646),
    @("Biomass [Solid]", 2030, $null, $null, 0.0002970049348188396, $null, $null, $null, $null, $null, $null),
    @("Renewable Energy Carrier", 2030, $null, $null, 0.0001076404163676046, $null, $null, $null, $null, $null, $null),
    @("Overall Demand", 2030, $null, $null, 0.0004641652259788762, 0.01828014017674983, [double]"4.731207815303437e-05", 0.003176425757448155, 0.01198822545005427, 0.0002508564729736633, 0.02350081407665946),
    @("Hydrogen", 2040, $null, $null, $null, 0.0004366967881524865, $null, [double]"9.584626307956239e-09", [double]"6.881349660176644e-05", $null, $null),
    @("Methanol", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Ammonia", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Synthetic Gases", 2040, $null, $null, $null, [double]"1.594077719567545e-10", $null, $null, [double]"2.065273219980239e-11", $null, $null),
    @("Biogenic Gases", 2040, $null, $null, 0.0002367686394466318, [double]"3.37213248451954e-05", $null, $null, [double]"1.41983412411986e-05", $null, $null),
    @("Fossil Gases", 2040, $null, $null, $null, 0.0001890745293606341, $null, $null, [double]"3.621938182632187e-05", $null, $null),
    @("Synthetic Liquids", 2040, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Biogenic Liquids", 2040, $null, $null, $null, 0.0007047696184429826, [double]"9.589618898468852e-06", 0.0003813658689076, 0.0010709142162089, [double]"4.320853681627712e-05", 0.002463983516709142),
    @("Fossil Liquids", 2040, $null, $null, $null, 0.00436006402907351, [double]"4.453898550539274e-05", 0.0026994173542879, 0.0045305831688898, 0.0001910498510613, 0.02068552430275745),
    @("Biomass [Solid]", 2040, $null, $null, 0.000295424140406074, $null, $null, $null, $null, $null, $null),
    @("Renewable Energy Carrier", 2040, $null, $null, 0.0004265001570010997, $null, $null, $null, $null, $null, $null),
    @("Overall Demand", 2040, $null, $null, 0.0009586929368538056, 0.005724326449282582, [double]"5.412860440386159e-05", 0.003080792807821808, 0.00572072862542072, 0.0002342583878775771, 0.02314950781946659),
    @("Hydrogen", 2050, $null, $null, $null, 0.0006056041171712137, $null, [double]"1.624516956572593e-08", 0.0001092982469662297, $null, $null),
    @("Methanol", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Ammonia", 2050, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @("Synthetic Gases", 2050, $null, $null, $null, [double]"1.448577107678552e-09", $null, $null, [double]"4.608565549328089e-10", $null, $null),
    @("Biogenic Gases", 2050, $null, $null, 0.0005637625130154859, [double]"5.71834062520828e-06", $null, $null, [double]"4.148322748065358e-06", $null, $null),
    @("Fossil Gases", 2050, $null, $null, $null, [double]"1.157307830130374e-05", $null, $null, [double]"1.387990866727912e-05", $null, $null),
    @("Synthetic Liquids", 2050, $null, $null, $null, [double]"3.021982398372452e-12", [double]"3.531338325264679e-13", [double]"1.469299066349851e-11", [double]"1.824420290111216e-11", [double]"2.959277376230371e-13", [double]"1.688469752281557e-10"),
    @("Biogenic Liquids", 2050, $null, $null, $null, [double]"6.371611843765214e-05", [double]"1.713496530849602e-05", 0.0004995071219427, 0.0002751207677780113, [double]"5.557302063796528e-05", 0.003510127756655649),
    @("Fossil Liquids", 2050, $null, $null, $null, 0.0002032153415099435, [double]"4.017227868407744e-05", 0.0024480240244066, 0.0008140281851726, 0.0001644463307791, 0.01928324458791736),
    @("Biomass [Solid]", 2050, $null, $null, 0.0002818113349595355, $null, $null, $null, $null, $null, $null),
    @("Renewable Energy Carrier", 2050, $null, $null, 0.0009826922911579049, $null, $null, $null, $null, $null, $null),
    @("Overall Demand", 2050, $null, $null, 0.001828266139132926, 0.0008898284476444115, [double]"5.730724434570729e-05", 0.002947547406211856, 0.001216475910432943, 0.000220019351712993, 0.02279337251341998)
)

$startRow = 7
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $val = $rowData[$c - 1]
        if ($null -eq $val) {
            $ws.Cells.Item($r, $c).Value = ""
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

Write-Host "Applied FuelGroup rows 7-37 (sheet dimension now A1:K37)"
